$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("I7").Value = "sd"
$ws.Range("J7").Value = "Statement-non-opinion"

# Row 18
$ws.Range("I18").Value = "sv"
$ws.Range("J18").Value = "Statement-opinion"

# Row 22
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"

# Row 23
$ws.Range("I23").Value = "sd"
$ws.Range("J23").Value = "Statement-non-opinion"

# Row 42
$ws.Range("I42").Value = "sv"
$ws.Range("J42").Value = "Statement-opinion"

# Row 51
$ws.Range("I51").Value = "sd"
$ws.Range("J51").Value = "Statement-non-opinion"

# Row 61
$ws.Range("I61").Value = "sv"
$ws.Range("J61").Value = "Statement-opinion"

# Row 66
$ws.Range("I66").Value = "sv"
$ws.Range("J66").Value = "Statement-opinion"

$wb.Save()
